# Recitation_12.pptx edit:
#  - slide 14: "0/1" -> "0|1" in the two outer cells of the allocation table
#  - slide 28: fix "Insert Block 4_5 from the Free List" -> "...to the Free List"
#  - slide 32: "0/1" -> "0|1" in the two outer cells of the allocation table,
#               and move the "Picture 29" image to the end of the shape order
#               (same embedded picture, new position/size) to match the
#               re-inserted <p:pic> that now sits right before </p:spTree>.

$p = $ppt.ActivePresentation

# ---- Slide 14 : table "0/1" -> "0|1" -----------------------------------
$s14 = $p.Slides.Item(14)
$tbl14 = $null
for ($i = 1; $i -le $s14.Shapes.Count; $i++) {
    $sh = $s14.Shapes.Item($i)
    if ($sh.HasTable) {
        $tbl14 = $sh.Table
        break
    }
}
$tbl14.Cell(1, 1).Shape.TextFrame.TextRange.Text = "0|1"
$tbl14.Cell(1, $tbl14.Columns.Count).Shape.TextFrame.TextRange.Text = "0|1"

# ---- Slide 28 : "from the Free List" -> "to the Free List" -------------
$s28 = $p.Slides.Item(28)
for ($i = 1; $i -le $s28.Shapes.Count; $i++) {
    $sh = $s28.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "Step 2: Insert Block 4_5 from the Free List") {
            $runCount = $sh.TextFrame.TextRange.Runs().Count
            for ($r = 1; $r -le $runCount; $r++) {
                $run = $sh.TextFrame.TextRange.Runs($r)
                if ($run.Text -eq " Insert Block 4_5 from the Free List") {
                    $run.Text = " Insert Block 4_5 to the Free List"
                }
            }
        }
    }
}

# ---- Slide 32 : table "0/1" -> "0|1" ------------------------------------
$s32 = $p.Slides.Item(32)
$tbl32 = $null
for ($i = 1; $i -le $s32.Shapes.Count; $i++) {
    $sh = $s32.Shapes.Item($i)
    if ($sh.HasTable) {
        $tbl32 = $sh.Table
        break
    }
}
$tbl32.Cell(1, 1).Shape.TextFrame.TextRange.Text = "0|1"
$tbl32.Cell(1, $tbl32.Columns.Count).Shape.TextFrame.TextRange.Text = "0|1"

# ---- Slide 32 : reposition/reorder the "Picture 29" image --------------
$pic = $null
for ($i = 1; $i -le $s32.Shapes.Count; $i++) {
    $sh = $s32.Shapes.Item($i)
    if ($sh.Name -eq "Picture 29") {
        $pic = $sh
        break
    }
}

if ($pic -ne $null) {
    # Bring it to the front of the z-order, i.e. make it the very last
    # shape in the slide's shape tree (after the footer placeholder).
    $pic.ZOrder(0)
    $pic.Name = "Picture 32"
    $pic.Left = 696.393952007874
    $pic.Top = 28.848426196850394
    $pic.Width = 262.2227639055118
    $pic.Height = 42.60370078740157
}
